# Update the LR-pairs data: target clusters now include "FAPs" (following Dr Hou advice).
# All per-row statistics are recomputed for the expanded cluster set, rows 2-9 are updated
# in place and two new rows (10-11) are appended so each sending cluster (M1, M2) now has
# five target clusters: ECs, FAPs, M1, M2, sCs.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: M1 -> ECs
$ws.Range("A2").Value = "M1"
$ws.Range("B2").Value = "Cd22"
$ws.Range("C2").Value = "Ptprc"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 5.348457333333333
$ws.Range("H2").Value = 16.045372
$ws.Range("I2").Value = 0.4996811083875221
$ws.Range("J2").Value = 0.499681108387522
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.7744373333333333
$ws.Range("N2").Value = 2.323312
$ws.Range("O2").Value = 0.001996953570362765
$ws.Range("P2").Value = 0.001996953570362765
$ws.Range("Q2").Value = 4.142045034673777
$ws.Range("R2").Value = 37.278405312064
$ws.Range("S2").Value = 0.0009978399734372858
$ws.Range("T2").Value = 0.000997839973437286

# Row 3: M1 -> FAPs
$ws.Range("A3").Value = "M1"
$ws.Range("B3").Value = "Cd22"
$ws.Range("C3").Value = "Ptprc"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 5.348457333333333
$ws.Range("H3").Value = 16.045372
$ws.Range("I3").Value = 0.4996811083875221
$ws.Range("J3").Value = 0.499681108387522
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.2864023333333334
$ws.Range("N3").Value = 0.8592070000000001
$ws.Range("O3").Value = 0.0007385131598040558
$ws.Range("P3").Value = 0.0007385131598040558
$ws.Range("Q3").Value = 1.531810660000445
$ws.Range("R3").Value = 13.786295940004
$ws.Range("S3").Value = 0.0003690210742496618
$ws.Range("T3").Value = 0.0003690210742496618

# Row 4: M1 -> M1
$ws.Range("A4").Value = "M1"
$ws.Range("B4").Value = "Cd22"
$ws.Range("C4").Value = "Ptprc"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 5.348457333333333
$ws.Range("H4").Value = 16.045372
$ws.Range("I4").Value = 0.4996811083875221
$ws.Range("J4").Value = 0.499681108387522
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 127.3845316666667
$ws.Range("N4").Value = 382.153595
$ws.Range("O4").Value = 0.3284720200998472
$ws.Range("P4").Value = 0.3284720200998472
$ws.Range("Q4").Value = 681.3107325458155
$ws.Range("R4").Value = 6131.79659291234
$ws.Range("S4").Value = 0.1641312630777801
$ws.Range("T4").Value = 0.16413126307778

# Row 5: M1 -> M2
$ws.Range("A5").Value = "M1"
$ws.Range("B5").Value = "Cd22"
$ws.Range("C5").Value = "Ptprc"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 5.348457333333333
$ws.Range("H5").Value = 16.045372
$ws.Range("I5").Value = 0.4996811083875221
$ws.Range("J5").Value = 0.499681108387522
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 259.1857756666666
$ws.Range("N5").Value = 777.557327
$ws.Range("O5").Value = 0.668332914526494
$ws.Range("P5").Value = 0.6683329145264941
$ws.Range("Q5").Value = 1386.244062560071
$ws.Range("R5").Value = 12476.19656304064
$ws.Range("S5").Value = 0.3339533315024616
$ws.Range("T5").Value = 0.3339533315024616

# Row 6: M1 -> sCs
$ws.Range("A6").Value = "M1"
$ws.Range("B6").Value = "Cd22"
$ws.Range("C6").Value = "Ptprc"
$ws.Range("D6").Value = "sCs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 5.348457333333333
$ws.Range("H6").Value = 16.045372
$ws.Range("I6").Value = 0.4996811083875221
$ws.Range("J6").Value = 0.499681108387522
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.1782366666666667
$ws.Range("N6").Value = 0.53471
$ws.Range("O6").Value = 0.0004595986434919951
$ws.Range("P6").Value = 0.0004595986434919951
$ws.Range("Q6").Value = 0.9532912069022222
$ws.Range("R6").Value = 8.579620862120001
$ws.Range("S6").Value = 0.0002296527595934817
$ws.Range("T6").Value = 0.0002296527595934817

# Row 7: M2 -> ECs
$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "Cd22"
$ws.Range("C7").Value = "Ptprc"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 5.355284000000001
$ws.Range("H7").Value = 16.065852
$ws.Range("I7").Value = 0.500318891612478
$ws.Range("J7").Value = 0.5003188916124779
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.7744373333333333
$ws.Range("N7").Value = 2.323312
$ws.Range("O7").Value = 0.001996953570362765
$ws.Range("P7").Value = 0.001996953570362765
$ws.Range("Q7").Value = 4.147331860202668
$ws.Range("R7").Value = 37.32598674182401
$ws.Range("S7").Value = 0.0009991135969254791
$ws.Range("T7").Value = 0.0009991135969254791

# Row 8: M2 -> FAPs
$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "Cd22"
$ws.Range("C8").Value = "Ptprc"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 5.355284000000001
$ws.Range("H8").Value = 16.065852
$ws.Range("I8").Value = 0.500318891612478
$ws.Range("J8").Value = 0.5003188916124779
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.2864023333333334
$ws.Range("N8").Value = 0.8592070000000001
$ws.Range("O8").Value = 0.0007385131598040558
$ws.Range("P8").Value = 0.0007385131598040558
$ws.Range("Q8").Value = 1.533765833262667
$ws.Range("R8").Value = 13.803892499364
$ws.Range("S8").Value = 0.000369492085554394
$ws.Range("T8").Value = 0.000369492085554394

# Row 9: M2 -> M1
$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "Cd22"
$ws.Range("C9").Value = "Ptprc"
$ws.Range("D9").Value = "M1"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 5.355284000000001
$ws.Range("H9").Value = 16.065852
$ws.Range("I9").Value = 0.500318891612478
$ws.Range("J9").Value = 0.5003188916124779
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 127.3845316666667
$ws.Range("N9").Value = 382.153595
$ws.Range("O9").Value = 0.3284720200998472
$ws.Range("P9").Value = 0.3284720200998472
$ws.Range("Q9").Value = 682.1803442819935
$ws.Range("R9").Value = 6139.623098537941
$ws.Range("S9").Value = 0.1643407570220671
$ws.Range("T9").Value = 0.1643407570220671

# Row 10: M2 -> M2
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Cd22"
$ws.Range("C10").Value = "Ptprc"
$ws.Range("D10").Value = "M2"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 5.355284000000001
$ws.Range("H10").Value = 16.065852
$ws.Range("I10").Value = 0.500318891612478
$ws.Range("J10").Value = 0.5003188916124779
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 259.1857756666666
$ws.Range("N10").Value = 777.557327
$ws.Range("O10").Value = 0.668332914526494
$ws.Range("P10").Value = 0.6683329145264941
$ws.Range("Q10").Value = 1388.013437455289
$ws.Range("R10").Value = 12492.12093709761
$ws.Range("S10").Value = 0.3343795830240325
$ws.Range("T10").Value = 0.3343795830240325

# Row 11: M2 -> sCs
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Cd22"
$ws.Range("C11").Value = "Ptprc"
$ws.Range("D11").Value = "sCs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 5.355284000000001
$ws.Range("H11").Value = 16.065852
$ws.Range("I11").Value = 0.500318891612478
$ws.Range("J11").Value = 0.5003188916124779
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.1782366666666667
$ws.Range("N11").Value = 0.53471
$ws.Range("O11").Value = 0.0004595986434919951
$ws.Range("P11").Value = 0.0004595986434919951
$ws.Range("Q11").Value = 0.9545079692133336
$ws.Range("R11").Value = 8.590571722920002
$ws.Range("S11").Value = 0.0002299458838985134
$ws.Range("T11").Value = 0.0002299458838985134

